$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PEBCOM")

# Fill in the two previously-empty "OT" (column E) cells.
$ws.Range("E2").Value = "Pendiente ADM"
$ws.Range("E64").Value = "Pendiente ADM"

# --- Append two new data rows (84 and 85) -------------------------------
# Columns A,B,C,D,E,F,G,H,J,K,L,O,P are plain text in this sheet (even when
# their contents look numeric, e.g. "6571" or "8/6/2025"); only I, M and N
# hold real numbers. Force the text columns to Text format BEFORE writing
# the values so Excel doesn't auto-convert things like dates or numeric
# looking strings, then clear the formatting again afterwards so the new
# cells end up unstyled, same as the rest of the data rows.

$textCols = "A", "B", "C", "D", "E", "F", "G", "H", "J", "K", "L", "O", "P"

foreach ($col in $textCols) {
    $ws.Range($col + "84").NumberFormat = "@"
    $ws.Range($col + "85").NumberFormat = "@"
}

# New row 84
$ws.Range("A84").Value = "6571"
$ws.Range("B84").Value = "8/6/2025"
$ws.Range("C84").Value = "BERUTI 2496"
$ws.Range("D84").Value = "2"
$ws.Range("E84").Value = "808733917"
$ws.Range("F84").Value = "PEBCOM"
$ws.Range("G84").Value = "Pendiente"
$ws.Range("H84").Value = "Aplomar"
$ws.Range("I84").Value = 1
$ws.Range("J84").Value = "Aplomo"
$ws.Range("K84").Value = "Sin equipos"
$ws.Range("L84").Value = "Pasante"
$ws.Range("M84").Value = -58.401374
$ws.Range("N84").Value = -34.592623
$ws.Range("O84").Value = "Recoleta"
$ws.Range("P84").Value = "Capital Sur"

# New row 85
$ws.Range("A85").Value = "6572"
$ws.Range("B85").Value = "8/6/2025"
$ws.Range("C85").Value = "MEXICO 2639"
$ws.Range("D85").Value = "3"
$ws.Range("E85").Value = "808733920"
$ws.Range("F85").Value = "PEBCOM"
$ws.Range("G85").Value = "Pendiente"
$ws.Range("H85").Value = "Picada"
$ws.Range("I85").Value = 1
$ws.Range("J85").Value = "Cambio"
$ws.Range("K85").Value = "Sin equipos"
$ws.Range("L85").Value = "Pasante"
$ws.Range("M85").Value = -58.403444
$ws.Range("N85").Value = -34.61685
$ws.Range("O85").Value = "Almagro"
$ws.Range("P85").Value = "Capital Sur"

# Strip the temporary Text number-format again so the new rows stay
# unstyled (like every other data row in the sheet).
$ws.Range("A84:P85").ClearFormats()
